# Update "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) timestamps on row 3
# (the e900ff5d-... entry) for both the zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-17 07:37:07"
$wsZhCn.Range("G3").Value = "2016-01-17 07:37:50"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-17 07:37:17"
$wsDeDe.Range("G3").Value = "2016-01-17 07:38:08"
